# Apply updated evaluation metrics to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": row 2 (model = isolation_forest/augmented/magnitude_5/split_2/test_50_50)
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.5
$wsSummary.Range("C2").Value = 0.5
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.6666666666666666
$wsSummary.Range("F2").Value = 0.8333333333333334
$wsSummary.Range("G2").Value = 0.9629629629629629
$wsSummary.Range("H2").Value = 0.7865010730968313
$wsSummary.Range("I2").Value = 534
$wsSummary.Range("J2").Value = 534
$wsSummary.Range("K2").Value = 0
$wsSummary.Range("L2").Value = 0

# --- Sheet "Classification Report"
$wsReport = $wb.Worksheets.Item("Classification Report")

# Row 2: class "0"
$wsReport.Range("B2").Value = 0
$wsReport.Range("C2").Value = 0
$wsReport.Range("D2").Value = 0

# Row 3: class "1"
$wsReport.Range("B3").Value = 0.5
$wsReport.Range("C3").Value = 1
$wsReport.Range("D3").Value = 0.6666666666666666

# Row 4: accuracy
$wsReport.Range("B4").Value = 0.5
$wsReport.Range("C4").Value = 0.5
$wsReport.Range("D4").Value = 0.5
$wsReport.Range("E4").Value = 0.5

# Row 5: macro avg
$wsReport.Range("B5").Value = 0.25
$wsReport.Range("C5").Value = 0.5
$wsReport.Range("D5").Value = 0.3333333333333333

# Row 6: weighted avg
$wsReport.Range("B6").Value = 0.25
$wsReport.Range("C6").Value = 0.5
$wsReport.Range("D6").Value = 0.3333333333333333

# --- Sheet "Confusion Matrix"
$wsConf = $wb.Worksheets.Item("Confusion Matrix")

# Row 2: Actual 0
$wsConf.Range("B2").Value = 0
$wsConf.Range("C2").Value = 534

# Row 3: Actual 1
$wsConf.Range("B3").Value = 0
$wsConf.Range("C3").Value = 534
